# Updates generated-output sheet at 456a3b4.
# - Sheet "展览": a new event (VWonderland) is prepended at row 25, pushing the
#   "coser动漫展" and "《碧蓝航线》..." rows down one slot (with small view-count
#   bumps), and the old "创世次元..." row's content is dropped entirely.
#   Every other existing row across all four sheets just gets its column-F
#   (view/click count) bumped slightly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Row 25: brand-new event inserted ahead of the old row-25 event.
$ws1.Range("C25").Value = "上海·VWonderland虚拟主播线下见面会"
$ws1.Range("D25").Value = "翔殷路1099号 合生汇"
$ws1.Range("E25").Value = "2024.08.24 10:00-08.24 21:00"
$ws1.Range("F25").Value = 0
$ws1.Range("G25").Value = 60
$ws1.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=90693"
$ws1.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202408/FZ9CsGO81723560782092.png"

# Row 26: now holds what used to be row 25 ("coser动漫展"), view count bumped.
$ws1.Range("C26").Value = "上海·coser动漫展"
$ws1.Range("D26").Value = "海潮路133号B1 JUMP工坊"
$ws1.Range("E26").Value = "2024.08.24 10:00-08.25 17:00"
$ws1.Range("F26").Value = 1705
$ws1.Range("G26").Value = 60
$ws1.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=87347"
$ws1.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202406/i6vAgX8I1719311206769.jpeg"

# Row 27: now holds what used to be row 26 ("《碧蓝航线》..."), view count bumped.
# The old row-27 event ("创世次元...") is removed from the sheet entirely.
$ws1.Range("C27").Value = "上海·《碧蓝航线》 2024港区盛夏清凉节"
$ws1.Range("D27").Value = "龙腾大道2350号 西岸穹顶艺术中心"
$ws1.Range("E27").Value = "2024.08.24 10:00-08.25 18:00"
$ws1.Range("F27").Value = 1689
$ws1.Range("G27").Value = "已售罄"
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=89864"
$ws1.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202408/ulbsTGKK1723249007529.jpeg"

# Remaining rows on this sheet: only the column-F view/click counter changes.
$sheet1Updates = @(
    @{Cell="F4"; Value=654},
    @{Cell="F5"; Value=2939},
    @{Cell="F6"; Value=21},
    @{Cell="F10"; Value=6891},
    @{Cell="F11"; Value=39},
    @{Cell="F12"; Value=70},
    @{Cell="F13"; Value=352},
    @{Cell="F14"; Value=601},
    @{Cell="F15"; Value=1490},
    @{Cell="F16"; Value=1112},
    @{Cell="F17"; Value=2231},
    @{Cell="F18"; Value=1478},
    @{Cell="F20"; Value=111},
    @{Cell="F21"; Value=1108},
    @{Cell="F22"; Value=122},
    @{Cell="F23"; Value=176},
    @{Cell="F24"; Value=337},
    @{Cell="F28"; Value=1030},
    @{Cell="F29"; Value=36},
    @{Cell="F31"; Value=1219},
    @{Cell="F33"; Value=582},
    @{Cell="F34"; Value=30},
    @{Cell="F36"; Value=423},
    @{Cell="F37"; Value=11},
    @{Cell="F38"; Value=2461},
    @{Cell="F39"; Value=2716},
    @{Cell="F40"; Value=72},
    @{Cell="F45"; Value=317},
    @{Cell="F46"; Value=121},
    @{Cell="F47"; Value=167},
    @{Cell="F48"; Value=153}
)
foreach ($u in $sheet1Updates) {
    $ws1.Range($u.Cell).Value = $u.Value
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @(
    @{Cell="F13"; Value=4},
    @{Cell="F14"; Value=60},
    @{Cell="F15"; Value=57},
    @{Cell="F17"; Value=165},
    @{Cell="F20"; Value=48},
    @{Cell="F23"; Value=468}
)
foreach ($u in $sheet2Updates) {
    $ws2.Range($u.Cell).Value = $u.Value
}

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @(
    @{Cell="F4"; Value=539},
    @{Cell="F6"; Value=1675},
    @{Cell="F7"; Value=1849},
    @{Cell="F8"; Value=2722},
    @{Cell="F9"; Value=1006},
    @{Cell="F10"; Value=920},
    @{Cell="F12"; Value=260},
    @{Cell="F13"; Value=1442},
    @{Cell="F14"; Value=7324}
)
foreach ($u in $sheet3Updates) {
    $ws3.Range($u.Cell).Value = $u.Value
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All Types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @(
    @{Cell="F3"; Value=539},
    @{Cell="F4"; Value=654},
    @{Cell="F5"; Value=2939},
    @{Cell="F6"; Value=1675},
    @{Cell="F8"; Value=2722},
    @{Cell="F9"; Value=6891},
    @{Cell="F10"; Value=1006},
    @{Cell="F11"; Value=39},
    @{Cell="F12"; Value=352},
    @{Cell="F14"; Value=260},
    @{Cell="F15"; Value=1112},
    @{Cell="F16"; Value=2231},
    @{Cell="F17"; Value=1478},
    @{Cell="F18"; Value=111},
    @{Cell="F20"; Value=1108},
    @{Cell="F21"; Value=60},
    @{Cell="F22"; Value=1705},
    @{Cell="F23"; Value=165},
    @{Cell="F24"; Value=36},
    @{Cell="F26"; Value=1219},
    @{Cell="F29"; Value=582},
    @{Cell="F30"; Value=30},
    @{Cell="F31"; Value=49},
    @{Cell="F34"; Value=468},
    @{Cell="F35"; Value=423},
    @{Cell="F37"; Value=11},
    @{Cell="F38"; Value=2461},
    @{Cell="F39"; Value=2716},
    @{Cell="F40"; Value=72},
    @{Cell="F44"; Value=317},
    @{Cell="F45"; Value=121},
    @{Cell="F46"; Value=167}
)
foreach ($u in $sheet4Updates) {
    $ws4.Range($u.Cell).Value = $u.Value
}
